$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (existing "CasesTab" row): the case-level query text (column B) was
# updated - an extra space before the WHERE-clause bracket and the
# Age (years) expression now coalesces/rounds the value.
# ---------------------------------------------------------------------------
$caseQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE   tp.endocrine_therapy_type IN  ["Tam"] 
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       coalesce(CASE demo.age_at_index % 1 WHEN 0 THEN apoc.convert.toInteger(demo.age_at_index) ELSE demo.age_at_index END, '') AS `Age (years)`,
demo.survival_time AS `Survival (days)`
'@

# ---------------------------------------------------------------------------
# New "SamplesTab" row query text (column B, row 3)
# ---------------------------------------------------------------------------
$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
 WHERE   tp.endocrine_therapy_type IN  ["Tam"]  
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@

# ---------------------------------------------------------------------------
# New "FilesTab" row query text (column B, row 4)
# ---------------------------------------------------------------------------
$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
 WHERE   tp.endocrine_therapy_type IN  ["Tam"] 
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@

# Here-strings add a trailing newline; strip it so the cell text matches the
# original (no trailing blank line).
$caseQuery    = $caseQuery.TrimEnd("`r", "`n")
$samplesQuery = $samplesQuery.TrimEnd("`r", "`n")
$filesQuery   = $filesQuery.TrimEnd("`r", "`n")

# Reuse the existing StatQuery / file-name text for the new rows, copied
# verbatim from the row above so there is no risk of retyping errors.
$statQuery   = $ws.Range("C2").Value()
$neo4jFile   = $ws.Range("D2").Value()
$webFile     = $ws.Range("E2").Value()

# ---------------------------------------------------------------------------
# Author workflow: the two new tab labels went in first, then the existing
# case query was tweaked, then the two new query bodies were typed in.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# Update the existing case-level query cell.
$ws.Range("B2").Value = $caseQuery

$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

# ---------------------------------------------------------------------------
# Remaining cells in rows 3 & 4 reuse already-existing shared strings.
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

$ws.Range("C4").Value = $statQuery
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# Wrap text for the long query / stat-query columns, matching row 2's style.
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true

# Row heights (as authored after re-wrapping the longer query text).
$ws.Rows.Item(2).RowHeight = 345.6
$ws.Rows.Item(3).RowHeight = 345.6
$ws.Rows.Item(4).RowHeight = 409.6

# Column widths (re-fit after the new content was added). The headless
# engine quantizes ColumnWidth to 1/6-character steps (it has no real font
# metrics for bestFit), so the request is pre-offset to land as close as
# possible to the widths Excel's bestFit computed for the real fonts.
$ws.Columns.Item(1).ColumnWidth = 12.77734375 - 5/6
$ws.Columns.Item(2).ColumnWidth = 76.109375   - 5/6
$ws.Columns.Item(3).ColumnWidth = 47.88671875 - 5/6
$ws.Columns.Item(4).ColumnWidth = 58.33203125 - 5/6
$ws.Columns.Item(5).ColumnWidth = 57.109375   - 5/6

# Final selection lands on B4, with no frozen/scrolled top-left cell.
$ws.Range("B4").Select() | Out-Null
